$wb = $excel.ActiveWorkbook

# Sheet "展览" (展览信息表)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 450
$ws.Range("F9").Value = 275
$ws.Range("F13").Value = 10
$ws.Range("F15").Value = 287
$ws.Range("F16").Value = 82
$ws.Range("F22").Value = 1639
$ws.Range("F23").Value = 396
$ws.Range("F26").Value = 1199
$ws.Range("F30").Value = 1587
$ws.Range("F33").Value = 621
$ws.Range("F34").Value = 853
$ws.Range("F35").Value = 1720
$ws.Range("F37").Value = 1750
$ws.Range("F40").Value = 829
$ws.Range("F41").Value = 30
$ws.Range("F44").Value = 983
$ws.Range("F46").Value = 3310

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 450
$ws.Range("F9").Value = 275
$ws.Range("F14").Value = 10
$ws.Range("F16").Value = 287
$ws.Range("F21").Value = 1639
$ws.Range("F23").Value = 396
$ws.Range("F27").Value = 1199
$ws.Range("F29").Value = 1587
$ws.Range("F36").Value = 621
$ws.Range("F37").Value = 1720
$ws.Range("F41").Value = 1750
$ws.Range("F42").Value = 829
$ws.Range("F45").Value = 983
$ws.Range("F48").Value = 3310
